$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column F "time_taken": copy header style from E1, then set values ---
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$ws.Range("F2").Value = "2021-10-05 13:42:36.929809"
$ws.Range("F3").Value = "2021-10-05 13:42:36.929823"
$ws.Range("F4").Value = "2021-10-05 13:42:36.929827"
$ws.Range("F5").Value = "2021-10-05 13:42:36.929829"
$ws.Range("F6").Value = "2021-10-05 13:42:36.929832"
$ws.Range("F7").Value = "2021-10-05 13:42:36.929835"
$ws.Range("F8").Value = "2021-10-05 13:42:36.929837"
$ws.Range("F9").Value = "2021-10-05 13:42:36.929840"
$ws.Range("F10").Value = "2021-10-05 13:42:36.929842"
$ws.Range("F11").Value = "2021-10-05 13:42:36.929845"
$ws.Range("F12").Value = "2021-10-05 13:42:36.929847"
$ws.Range("F13").Value = "2021-10-05 13:42:36.929850"
$ws.Range("F14").Value = "2021-10-05 13:42:36.929852"
$ws.Range("F15").Value = "2021-10-05 13:42:36.929855"
$ws.Range("F16").Value = "2021-10-05 13:42:36.929857"
$ws.Range("F17").Value = "2021-10-05 13:42:36.929859"
$ws.Range("F18").Value = "2021-10-05 13:42:36.929862"
$ws.Range("F19").Value = "2021-10-05 13:42:36.929865"
$ws.Range("F20").Value = "2021-10-05 13:42:36.929867"
$ws.Range("F21").Value = "2021-10-05 13:42:36.929870"
$ws.Range("F22").Value = "2021-10-05 13:42:36.929872"
$ws.Range("F23").Value = "2021-10-05 13:42:36.929875"
$ws.Range("F24").Value = "2021-10-05 13:42:36.929877"
$ws.Range("F25").Value = "2021-10-05 13:42:36.929880"
$ws.Range("F26").Value = "2021-10-05 13:42:36.929883"
$ws.Range("F27").Value = "2021-10-05 13:42:36.929885"
$ws.Range("F28").Value = "2021-10-05 13:42:36.929888"
$ws.Range("F29").Value = "2021-10-05 13:42:36.929890"
$ws.Range("F30").Value = "2021-10-05 13:42:36.929892"
$ws.Range("F31").Value = "2021-10-05 13:42:36.929895"
$ws.Range("F32").Value = "2021-10-05 13:42:36.929897"
$ws.Range("F33").Value = "2021-10-05 13:42:36.929900"
$ws.Range("F34").Value = "2021-10-05 13:42:36.929902"
$ws.Range("F35").Value = "2021-10-05 13:42:36.929905"
$ws.Range("F36").Value = "2021-10-05 13:42:36.929907"
$ws.Range("F37").Value = "2021-10-05 13:42:36.929910"
$ws.Range("F38").Value = "2021-10-05 13:42:36.929912"
$ws.Range("F39").Value = "2021-10-05 13:42:36.929915"
$ws.Range("F40").Value = "2021-10-05 13:42:36.929917"
$ws.Range("F41").Value = "2021-10-05 13:42:36.929919"
$ws.Range("F42").Value = "2021-10-05 13:42:36.929922"
$ws.Range("F43").Value = "2021-10-05 13:42:36.929925"
$ws.Range("F44").Value = "2021-10-05 13:42:36.929927"
$ws.Range("F45").Value = "2021-10-05 13:42:36.929930"
$ws.Range("F46").Value = "2021-10-05 13:42:36.929932"
$ws.Range("F47").Value = "2021-10-05 13:42:36.929935"
$ws.Range("F48").Value = "2021-10-05 13:42:36.929937"
$ws.Range("F49").Value = "2021-10-05 13:42:36.929940"
$ws.Range("F50").Value = "2021-10-05 13:42:36.929942"
$ws.Range("F51").Value = "2021-10-05 13:42:36.929950"
$ws.Range("F52").Value = "2021-10-05 13:42:36.929953"
$ws.Range("F53").Value = "2021-10-05 13:42:36.929956"
$ws.Range("F54").Value = "2021-10-05 13:42:36.929959"
$ws.Range("F55").Value = "2021-10-05 13:42:36.929961"
$ws.Range("F56").Value = "2021-10-05 13:42:36.929964"
$ws.Range("F57").Value = "2021-10-05 13:42:36.929966"
$ws.Range("F58").Value = "2021-10-05 13:42:36.929968"
$ws.Range("F59").Value = "2021-10-05 13:42:36.929971"
$ws.Range("F60").Value = "2021-10-05 13:42:36.929973"
$ws.Range("F61").Value = "2021-10-05 13:42:36.929976"
$ws.Range("F62").Value = "2021-10-05 13:42:36.929978"
$ws.Range("F63").Value = "2021-10-05 13:42:36.929981"
$ws.Range("F64").Value = "2021-10-05 13:42:36.929983"
$ws.Range("F65").Value = "2021-10-05 13:42:36.929985"
$ws.Range("F66").Value = "2021-10-05 13:42:36.929989"
$ws.Range("F67").Value = "2021-10-05 13:42:36.929992"
$ws.Range("F68").Value = "2021-10-05 13:42:36.929995"
$ws.Range("F69").Value = "2021-10-05 13:42:36.929997"
$ws.Range("F70").Value = "2021-10-05 13:42:36.929999"
$ws.Range("F71").Value = "2021-10-05 13:42:36.930002"
$ws.Range("F72").Value = "2021-10-05 13:42:36.930005"
$ws.Range("F73").Value = "2021-10-05 13:42:36.930007"
$ws.Range("F74").Value = "2021-10-05 13:42:36.930010"
$ws.Range("F75").Value = "2021-10-05 13:42:36.930012"
$ws.Range("F76").Value = "2021-10-05 13:42:36.930014"
$ws.Range("F77").Value = "2021-10-05 13:42:36.930017"
$ws.Range("F78").Value = "2021-10-05 13:42:36.930021"
$ws.Range("F79").Value = "2021-10-05 13:42:36.930024"
$ws.Range("F80").Value = "2021-10-05 13:42:36.930026"
$ws.Range("F81").Value = "2021-10-05 13:42:36.930029"
$ws.Range("F82").Value = "2021-10-05 13:42:36.930031"
$ws.Range("F83").Value = "2021-10-05 13:42:36.930034"
$ws.Range("F84").Value = "2021-10-05 13:42:36.930036"
$ws.Range("F85").Value = "2021-10-05 13:42:36.930039"
$ws.Range("F86").Value = "2021-10-05 13:42:36.930041"
$ws.Range("F87").Value = "2021-10-05 13:42:36.930043"
$ws.Range("F88").Value = "2021-10-05 13:42:36.930046"
$ws.Range("F89").Value = "2021-10-05 13:42:36.930048"
$ws.Range("F90").Value = "2021-10-05 13:42:36.930051"
$ws.Range("F91").Value = "2021-10-05 13:42:36.930053"
$ws.Range("F92").Value = "2021-10-05 13:42:36.930055"
$ws.Range("F93").Value = "2021-10-05 13:42:36.930058"
$ws.Range("F94").Value = "2021-10-05 13:42:36.930061"
$ws.Range("F95").Value = "2021-10-05 13:42:36.930064"
$ws.Range("F96").Value = "2021-10-05 13:42:36.930066"
$ws.Range("F97").Value = "2021-10-05 13:42:36.930069"
$ws.Range("F98").Value = "2021-10-05 13:42:36.930071"
$ws.Range("F99").Value = "2021-10-05 13:42:36.930074"
$ws.Range("F100").Value = "2021-10-05 13:42:36.930076"
$ws.Range("F101").Value = "2021-10-05 13:42:36.930078"
$ws.Range("F102").Value = "2021-10-05 13:42:36.930081"
$ws.Range("F103").Value = "2021-10-05 13:42:36.930083"
$ws.Range("F104").Value = "2021-10-05 13:42:36.930086"
$ws.Range("F105").Value = "2021-10-05 13:42:36.930088"
$ws.Range("F106").Value = "2021-10-05 13:42:36.930091"
$ws.Range("F107").Value = "2021-10-05 13:42:36.930093"
$ws.Range("F108").Value = "2021-10-05 13:42:36.930095"
$ws.Range("F109").Value = "2021-10-05 13:42:36.930098"
$ws.Range("F110").Value = "2021-10-05 13:42:36.930102"
$ws.Range("F111").Value = "2021-10-05 13:42:36.930105"
$ws.Range("F112").Value = "2021-10-05 13:42:36.930108"
$ws.Range("F113").Value = "2021-10-05 13:42:36.930110"
$ws.Range("F114").Value = "2021-10-05 13:42:36.930112"
$ws.Range("F115").Value = "2021-10-05 13:42:36.930115"
$ws.Range("F116").Value = "2021-10-05 13:42:36.930117"
$ws.Range("F117").Value = "2021-10-05 13:42:36.930120"
$ws.Range("F118").Value = "2021-10-05 13:42:36.930122"
$ws.Range("F119").Value = "2021-10-05 13:42:36.930125"
$ws.Range("F120").Value = "2021-10-05 13:42:36.930127"
$ws.Range("F121").Value = "2021-10-05 13:42:36.930129"
$ws.Range("F122").Value = "2021-10-05 13:42:36.930132"
$ws.Range("F123").Value = "2021-10-05 13:42:36.930134"
$ws.Range("F124").Value = "2021-10-05 13:42:36.930137"
$ws.Range("F125").Value = "2021-10-05 13:42:36.930139"
$ws.Range("F126").Value = "2021-10-05 13:42:36.930142"
$ws.Range("F127").Value = "2021-10-05 13:42:36.930144"
$ws.Range("F128").Value = "2021-10-05 13:42:36.930146"
$ws.Range("F129").Value = "2021-10-05 13:42:36.930149"
$ws.Range("F130").Value = "2021-10-05 13:42:36.930153"
$ws.Range("F131").Value = "2021-10-05 13:42:36.930156"
$ws.Range("F132").Value = "2021-10-05 13:42:36.930158"
$ws.Range("F133").Value = "2021-10-05 13:42:36.930161"
$ws.Range("F134").Value = "2021-10-05 13:42:36.930163"
$ws.Range("F135").Value = "2021-10-05 13:42:36.930166"
$ws.Range("F136").Value = "2021-10-05 13:42:36.930168"
$ws.Range("F137").Value = "2021-10-05 13:42:36.930170"
$ws.Range("F138").Value = "2021-10-05 13:42:36.930173"
$ws.Range("F139").Value = "2021-10-05 13:42:36.930175"
$ws.Range("F140").Value = "2021-10-05 13:42:36.930178"
$ws.Range("F141").Value = "2021-10-05 13:42:36.930180"
$ws.Range("F142").Value = "2021-10-05 13:42:36.930183"
$ws.Range("F143").Value = "2021-10-05 13:42:36.930185"
$ws.Range("F144").Value = "2021-10-05 13:42:36.930188"
$ws.Range("F145").Value = "2021-10-05 13:42:36.930190"
$ws.Range("F146").Value = "2021-10-05 13:42:36.930193"
$ws.Range("F147").Value = "2021-10-05 13:42:36.930195"
$ws.Range("F148").Value = "2021-10-05 13:42:36.930197"
$ws.Range("F149").Value = "2021-10-05 13:42:36.930200"
$ws.Range("F150").Value = "2021-10-05 13:42:36.930202"
$ws.Range("F151").Value = "2021-10-05 13:42:36.930205"
$ws.Range("F152").Value = "2021-10-05 13:42:36.930207"
$ws.Range("F153").Value = "2021-10-05 13:42:36.930210"
$ws.Range("F154").Value = "2021-10-05 13:42:36.930212"
$ws.Range("F155").Value = "2021-10-05 13:42:36.930215"
$ws.Range("F156").Value = "2021-10-05 13:42:36.930217"
$ws.Range("F157").Value = "2021-10-05 13:42:36.930220"
$ws.Range("F158").Value = "2021-10-05 13:42:36.930222"
$ws.Range("F159").Value = "2021-10-05 13:42:36.930225"
$ws.Range("F160").Value = "2021-10-05 13:42:36.930227"
$ws.Range("F161").Value = "2021-10-05 13:42:36.930229"
$ws.Range("F162").Value = "2021-10-05 13:42:36.930232"
$ws.Range("F163").Value = "2021-10-05 13:42:36.930234"
$ws.Range("F164").Value = "2021-10-05 13:42:36.930237"
$ws.Range("F165").Value = "2021-10-05 13:42:36.930239"
$ws.Range("F166").Value = "2021-10-05 13:42:36.930241"
$ws.Range("F167").Value = "2021-10-05 13:42:36.930244"
$ws.Range("F168").Value = "2021-10-05 13:42:36.930246"
$ws.Range("F169").Value = "2021-10-05 13:42:36.930249"
$ws.Range("F170").Value = "2021-10-05 13:42:36.930251"
$ws.Range("F171").Value = "2021-10-05 13:42:36.930253"
$ws.Range("F172").Value = "2021-10-05 13:42:36.930256"
$ws.Range("F173").Value = "2021-10-05 13:42:36.930258"
$ws.Range("F174").Value = "2021-10-05 13:42:36.930262"
$ws.Range("F175").Value = "2021-10-05 13:42:36.930265"
$ws.Range("F176").Value = "2021-10-05 13:42:36.930267"
$ws.Range("F177").Value = "2021-10-05 13:42:36.930270"
$ws.Range("F178").Value = "2021-10-05 13:42:36.930272"
$ws.Range("F179").Value = "2021-10-05 13:42:36.930275"
$ws.Range("F180").Value = "2021-10-05 13:42:36.930277"
$ws.Range("F181").Value = "2021-10-05 13:42:36.930280"
$ws.Range("F182").Value = "2021-10-05 13:42:36.930282"
$ws.Range("F183").Value = "2021-10-05 13:42:36.930284"
$ws.Range("F184").Value = "2021-10-05 13:42:36.930287"
$ws.Range("F185").Value = "2021-10-05 13:42:36.930289"
$ws.Range("F186").Value = "2021-10-05 13:42:36.930292"
$ws.Range("F187").Value = "2021-10-05 13:42:36.930294"
$ws.Range("F188").Value = "2021-10-05 13:42:36.930297"
$ws.Range("F189").Value = "2021-10-05 13:42:36.930299"
$ws.Range("F190").Value = "2021-10-05 13:42:36.930302"
$ws.Range("F191").Value = "2021-10-05 13:42:36.930304"
$ws.Range("F192").Value = "2021-10-05 13:42:36.930307"
$ws.Range("F193").Value = "2021-10-05 13:42:36.930309"
$ws.Range("F194").Value = "2021-10-05 13:42:36.930311"
$ws.Range("F195").Value = "2021-10-05 13:42:36.930314"
$ws.Range("F196").Value = "2021-10-05 13:42:36.930316"
$ws.Range("F197").Value = "2021-10-05 13:42:36.930319"
$ws.Range("F198").Value = "2021-10-05 13:42:36.930321"
$ws.Range("F199").Value = "2021-10-05 13:42:36.930323"
$ws.Range("F200").Value = "2021-10-05 13:42:36.930326"
$ws.Range("F201").Value = "2021-10-05 13:42:36.930329"
$ws.Range("F202").Value = "2021-10-05 13:42:36.930331"
$ws.Range("F203").Value = "2021-10-05 13:42:36.930334"
$ws.Range("F204").Value = "2021-10-05 13:42:36.930336"
$ws.Range("F205").Value = "2021-10-05 13:42:36.930339"
$ws.Range("F206").Value = "2021-10-05 13:42:36.930341"
$ws.Range("F207").Value = "2021-10-05 13:42:36.930344"
$ws.Range("F208").Value = "2021-10-05 13:42:36.930346"
$ws.Range("F209").Value = "2021-10-05 13:42:36.930349"
$ws.Range("F210").Value = "2021-10-05 13:42:36.930351"
$ws.Range("F211").Value = "2021-10-05 13:42:36.930353"
$ws.Range("F212").Value = "2021-10-05 13:42:36.930356"
$ws.Range("F213").Value = "2021-10-05 13:42:36.930358"
$ws.Range("F214").Value = "2021-10-05 13:42:36.930361"
$ws.Range("F215").Value = "2021-10-05 13:42:36.930363"
$ws.Range("F216").Value = "2021-10-05 13:42:36.930365"
$ws.Range("F217").Value = "2021-10-05 13:42:36.930368"
$ws.Range("F218").Value = "2021-10-05 13:42:36.930370"
$ws.Range("F219").Value = "2021-10-05 13:42:36.930373"
$ws.Range("F220").Value = "2021-10-05 13:42:36.930376"
$ws.Range("F221").Value = "2021-10-05 13:42:36.930378"
$ws.Range("F222").Value = "2021-10-05 13:42:36.930381"
$ws.Range("F223").Value = "2021-10-05 13:42:36.930383"
$ws.Range("F224").Value = "2021-10-05 13:42:36.930386"
$ws.Range("F225").Value = "2021-10-05 13:42:36.930388"
$ws.Range("F226").Value = "2021-10-05 13:42:36.930391"
$ws.Range("F227").Value = "2021-10-05 13:42:36.930393"
$ws.Range("F228").Value = "2021-10-05 13:42:36.930396"
$ws.Range("F229").Value = "2021-10-05 13:42:36.930398"

# --- Refresh geneSymbol (B) / geneName (C) for rows 157-199: row 157s data
# (GTF2E2) moves to the end (row 199), and rows 158-199 shift up by one ---
$ws.Range("B157").Value = "HIKESHI"
$ws.Range("C157").Value = "Hikeshi, heat shock protein nuclear import factor"
$ws.Range("B158").Value = "HSPD1"
$ws.Range("C158").Value = "heat shock protein family D (Hsp60) member 1"
$ws.Range("B159").Value = "ISCA1"
$ws.Range("C159").Value = "iron-sulfur cluster assembly 1"
$ws.Range("B160").Value = "ISCA2"
$ws.Range("C160").Value = "iron-sulfur cluster assembly 2"
$ws.Range("B161").Value = "KIAA1161"
$ws.Range("C161").Value = "myogenesis regulating glycosidase (putative)"
$ws.Range("B162").Value = "KIF5A"
$ws.Range("C162").Value = "kinesin family member 5A"
$ws.Range("B163").Value = "LIG3"
$ws.Range("C163").Value = "DNA ligase 3"
$ws.Range("B164").Value = "MRE11"
$ws.Range("C164").Value = "MRE11 homolog, double strand break repair nuclease"
$ws.Range("B165").Value = "MRPS16"
$ws.Range("C165").Value = "mitochondrial ribosomal protein S16"
$ws.Range("B166").Value = "NAXD"
$ws.Range("C166").Value = "NAD(P)HX dehydratase"
$ws.Range("B167").Value = "NAXE"
$ws.Range("C167").Value = "NAD(P)HX epimerase"
$ws.Range("B168").Value = "NFU1"
$ws.Range("C168").Value = "NFU1 iron-sulfur cluster scaffold"
$ws.Range("B169").Value = "NUP188"
$ws.Range("C169").Value = "nucleoporin 188"
$ws.Range("B170").Value = "PEX14"
$ws.Range("C170").Value = "peroxisomal biogenesis factor 14"
$ws.Range("B171").Value = "PEX19"
$ws.Range("C171").Value = "peroxisomal biogenesis factor 19"
$ws.Range("B172").Value = "PI4KA"
$ws.Range("C172").Value = "phosphatidylinositol 4-kinase alpha"
$ws.Range("B173").Value = "POLR3K"
$ws.Range("C173").Value = "RNA polymerase III subunit K"
$ws.Range("B174").Value = "PTEN"
$ws.Range("C174").Value = "phosphatase and tensin homolog"
$ws.Range("B175").Value = "RAB11B"
$ws.Range("C175").Value = "RAB11B, member RAS oncogene family"
$ws.Range("B176").Value = "RNU7-1"
$ws.Range("C176").Value = "RNA, U7 small nuclear 1"
$ws.Range("B177").Value = "RPIA"
$ws.Range("C177").Value = "ribose 5-phosphate isomerase A"
$ws.Range("B178").Value = "SCAF4"
$ws.Range("C178").Value = "SR-related CTD associated factor 4"
$ws.Range("B179").Value = "SDHA"
$ws.Range("C179").Value = "succinate dehydrogenase complex flavoprotein subunit A"
$ws.Range("B180").Value = "SNORD118"
$ws.Range("C180").Value = "small nucleolar RNA, C/D box 118"
$ws.Range("B181").Value = "SPART"
$ws.Range("C181").Value = "spartin"
$ws.Range("B182").Value = "SPG11"
$ws.Range("C182").Value = "SPG11, spatacsin vesicle trafficking associated"
$ws.Range("B183").Value = "STN1"
$ws.Range("C183").Value = "STN1, CST complex subunit"
$ws.Range("B184").Value = "TMEM106B"
$ws.Range("C184").Value = "transmembrane protein 106B"
$ws.Range("B185").Value = "TMEM63A"
$ws.Range("C185").Value = "transmembrane protein 63A"
$ws.Range("B186").Value = "TUFM"
$ws.Range("C186").Value = "Tu translation elongation factor, mitochondrial"
$ws.Range("B187").Value = "UFM1"
$ws.Range("C187").Value = "ubiquitin fold modifier 1"
$ws.Range("B188").Value = "VPS11"
$ws.Range("C188").Value = "VPS11, CORVET/HOPS core subunit"
$ws.Range("B189").Value = "WARS2"
$ws.Range("C189").Value = "tryptophanyl tRNA synthetase 2, mitochondrial"
$ws.Range("B190").Value = "ZFYVE26"
$ws.Range("C190").Value = "zinc finger FYVE-type containing 26"
$ws.Range("B191").Value = "ADGRG1"
$ws.Range("C191").Value = "adhesion G protein-coupled receptor G1"
$ws.Range("B192").Value = "ARX"
$ws.Range("C192").Value = "aristaless related homeobox"
$ws.Range("B193").Value = "ATP7A"
$ws.Range("C193").Value = "ATPase copper transporting alpha"
$ws.Range("B194").Value = "CYP7B1"
$ws.Range("C194").Value = "cytochrome P450 family 7 subfamily B member 1"
$ws.Range("B195").Value = "DCX"
$ws.Range("C195").Value = "doublecortin"
$ws.Range("B196").Value = "DDB1"
$ws.Range("C196").Value = "damage specific DNA binding protein 1"
$ws.Range("B197").Value = "EGR2"
$ws.Range("C197").Value = "early growth response 2"
$ws.Range("B198").Value = "FARSA"
$ws.Range("C198").Value = "phenylalanyl-tRNA synthetase alpha subunit"
$ws.Range("B199").Value = "GTF2E2"
$ws.Range("C199").Value = "general transcription factor IIE subunit 2"

# --- geneConfidence (D) at row 191 becomes 1 (ADGRG1 shifted in from row 192) ---
$ws.Range("D191").Value = "1"
